$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the device/preferences timezone-aware dateTool.format(...) helper instead of
# constructing a raw org.joda.time.DateTime / calling the bare from/to toString(),
# and switch the date format to "YYYY-MM-dd HH:mm:ss".

# Period value (row 6, column B) - combined from/to range
$ws.Range("B6").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", from, locale, timezone)+" - "+dateTool.format("YYYY-MM-dd HH:mm:ss", to, locale, timezone)}'

# Trip table header row (row 9): Start (A9) and End (C9) formulas
$ws.Range("A9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", trip.startTime, locale, timezone)}'
$ws.Range("C9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", trip.endTime, locale, timezone)}'

# Move the saved selection/active cell from D9 to B2
$ws.Range("B2").Select()
